$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 91.25
$ws.Range("I9").Value = 89.5
$ws.Range("K9").Value = 89.5
$ws.Range("M9").Value = 79.5

$ws.Range("H40").Value = 2666.3333
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2999.5
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2999.5
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -3349.5

$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H64").Value = 5357.25
$ws.Range("I64").Value = 3621.6667
$ws.Range("J64").Value = 6398.6
$ws.Range("K64").Value = 3621.6667
$ws.Range("L64").Value = 6398.6
$ws.Range("M64").Value = -3373.6667
$ws.Range("N64").Value = -6894.6

$ws.Range("H67").Value = 5357.25
$ws.Range("I67").Value = 3621.6667
$ws.Range("J67").Value = 6398.6
$ws.Range("K67").Value = 3621.6667
$ws.Range("L67").Value = 6398.6
$ws.Range("M67").Value = -2763.6667
$ws.Range("N67").Value = -8114.6

$ws.Range("H101").Value = 506
$ws.Range("I101").Value = 333
$ws.Range("J101").Value = 592.5
$ws.Range("K101").Value = 999
$ws.Range("L101").Value = 1777.5
$ws.Range("M101").Value = 623
$ws.Range("N101").Value = -5021.5

$ws.Range("H132").Value = 15294.421
$ws.Range("I132").Value = 18784.5
$ws.Range("J132").Value = 5522.2
$ws.Range("K132").Value = 56353.5
$ws.Range("L132").Value = 16566.6
$ws.Range("M132").Value = -53823.5
$ws.Range("N132").Value = -21626.6

$ws.Range("H137").Value = 14720481
$ws.Range("I137").Value = 31253286
$ws.Range("J137").Value = 24654.223
$ws.Range("K137").Value = 93759858
$ws.Range("L137").Value = 73962.66900000001
$ws.Range("M137").Value = -93757308
$ws.Range("N137").Value = -79062.66900000001

$ws.Range("H138").Value = 8435.59
$ws.Range("I138").Value = 10123.5
$ws.Range("K138").Value = 30370.5
$ws.Range("M138").Value = -25230.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 911803.1
$ws.Range("I32").Value = 973947.25
$ws.Range("J32").Value = 26249.75
$ws.Range("K32").Value = 973947.25
$ws.Range("L32").Value = 26249.75
$ws.Range("M32").Value = -973660.25
$ws.Range("N32").Value = -26823.75

$ws.Range("H74").Value = 673777.0600000001
$ws.Range("I74").Value = 779455.8
$ws.Range("J74").Value = 16220.556
$ws.Range("K74").Value = 779455.8
$ws.Range("L74").Value = 16220.556
$ws.Range("M74").Value = -778581.8
$ws.Range("N74").Value = -17968.556

$ws.Range("H77").Value = 673777.0600000001
$ws.Range("I77").Value = 779455.8
$ws.Range("J77").Value = 16220.556
$ws.Range("K77").Value = 3897279
$ws.Range("L77").Value = 81102.78
$ws.Range("M77").Value = -3892911
$ws.Range("N77").Value = -89838.78

$ws.Range("H102").Value = 2200.074
$ws.Range("I102").Value = 2207.8462
$ws.Range("J102").Value = 1998
$ws.Range("K102").Value = 2207.8462
$ws.Range("L102").Value = 1998
$ws.Range("M102").Value = -585.8462
$ws.Range("N102").Value = -5242

$ws.Range("H132").Value = 5329.9775
$ws.Range("I132").Value = 4205.65
$ws.Range("J132").Value = 6229.44
$ws.Range("K132").Value = 12616.95
$ws.Range("L132").Value = 18688.32
$ws.Range("M132").Value = -10086.95
$ws.Range("N132").Value = -23748.32

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 15941.7
$ws.Range("I99").Value = 18488.857
$ws.Range("J99").Value = 9998.333000000001
$ws.Range("K99").Value = 18488.857
$ws.Range("L99").Value = 9998.333000000001
$ws.Range("M99").Value = -16990.857
$ws.Range("N99").Value = -12994.333

$ws.Range("H134").Value = 3209601.5
$ws.Range("I134").Value = 4037.3
$ws.Range("J134").Value = 13894816
$ws.Range("K134").Value = 12111.9
$ws.Range("L134").Value = 41684448
$ws.Range("M134").Value = -9576.900000000001
$ws.Range("N134").Value = -41689518

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 494519.94
$ws.Range("I31").Value = 646132.1
$ws.Range("J31").Value = 4695.923
$ws.Range("K31").Value = 646132.1
$ws.Range("L31").Value = 4695.923
$ws.Range("M31").Value = -645837.1
$ws.Range("N31").Value = -5285.923

$ws.Range("H34").Value = 494519.94
$ws.Range("I34").Value = 646132.1
$ws.Range("J34").Value = 4695.923
$ws.Range("K34").Value = 646132.1
$ws.Range("L34").Value = 4695.923
$ws.Range("M34").Value = -645930.1
$ws.Range("N34").Value = -5099.923

$ws.Range("H132").Value = 2752.6316
$ws.Range("I132").Value = 2788.353
$ws.Range("J132").Value = 2449
$ws.Range("K132").Value = 8365.059000000001
$ws.Range("L132").Value = 7347
$ws.Range("M132").Value = -5835.059000000001
$ws.Range("N132").Value = -12407

$ws.Range("H134").Value = 4660.4
$ws.Range("I134").Value = 2191
$ws.Range("J134").Value = 5932.515
$ws.Range("K134").Value = 6573
$ws.Range("L134").Value = 17797.545
$ws.Range("M134").Value = -4038
$ws.Range("N134").Value = -22867.545

$ws.Range("H135").Value = 85999.336
$ws.Range("J135").Value = 85999.336
$ws.Range("L135").Value = 85999.336
$ws.Range("N135").Value = -96139.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 414.2857
$ws.Range("I8").Value = 414.2857
$ws.Range("K8").Value = 1242.8571
$ws.Range("M8").Value = -1103.8571

$ws.Range("H68").Value = 5610.0537
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 5639.327
$ws.Range("K68").Value = 12000
$ws.Range("L68").Value = 16917.981
$ws.Range("M68").Value = -11189
$ws.Range("N68").Value = -18539.981

$ws.Range("H71").Value = 5610.0537
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 5639.327
$ws.Range("K71").Value = 36000
$ws.Range("L71").Value = 50753.943
$ws.Range("M71").Value = -31944
$ws.Range("N71").Value = -58865.943

$ws.Range("H92").Value = 897.9167
$ws.Range("J92").Value = 1693.75
$ws.Range("L92").Value = 5081.25
$ws.Range("N92").Value = -7577.25

$ws.Range("H113").Value = 1604.125
$ws.Range("I113").Value = 758.75
$ws.Range("J113").Value = 2449.5
$ws.Range("K113").Value = 2276.25
$ws.Range("L113").Value = 7348.5
$ws.Range("M113").Value = -106.25
$ws.Range("N113").Value = -11688.5

$ws.Range("H136").Value = 9706.066000000001
$ws.Range("I136").Value = 5823.875
$ws.Range("J136").Value = 14142.857
$ws.Range("K136").Value = 17471.625
$ws.Range("L136").Value = 42428.571
$ws.Range("M136").Value = -12371.625
$ws.Range("N136").Value = -52628.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2402.875
$ws.Range("I113").Value = 2402.875
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2402.875
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -232.875
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 8763.529
$ws.Range("I132").Value = 11178.772
$ws.Range("J132").Value = 6931.276
$ws.Range("K132").Value = 33536.31600000001
$ws.Range("L132").Value = 20793.828
$ws.Range("M132").Value = -31006.31600000001
$ws.Range("N132").Value = -25853.828

$ws.Range("H135").Value = 88775.06
$ws.Range("J135").Value = 88775.06
$ws.Range("L135").Value = 88775.06
$ws.Range("N135").Value = -98915.06

$ws.Range("H136").Value = 59624.332
$ws.Range("J136").Value = 59624.332
$ws.Range("L136").Value = 178872.996
$ws.Range("N136").Value = -183972.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2381.7273
$ws.Range("I68").Value = 2399.9
$ws.Range("J68").Value = 2200
$ws.Range("K68").Value = 2399.9
$ws.Range("L68").Value = 2200
$ws.Range("M68").Value = -1650.9
$ws.Range("N68").Value = -3698

$ws.Range("H71").Value = 2381.7273
$ws.Range("I71").Value = 2399.9
$ws.Range("J71").Value = 2200
$ws.Range("K71").Value = 11999.5
$ws.Range("L71").Value = 11000
$ws.Range("M71").Value = -8255.5
$ws.Range("N71").Value = -18488

$ws.Range("H132").Value = 5559742.5
$ws.Range("I132").Value = 6948481
$ws.Range("J132").Value = 4786.6665
$ws.Range("K132").Value = 20845443
$ws.Range("L132").Value = 14359.9995
$ws.Range("M132").Value = -20842913
$ws.Range("N132").Value = -19419.9995

$ws.Range("H136").Value = 10872153
$ws.Range("I136").Value = 7814944
$ws.Range("J136").Value = 17860060
$ws.Range("K136").Value = 23444832
$ws.Range("L136").Value = 53580180
$ws.Range("M136").Value = -23442282
$ws.Range("N136").Value = -53585280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 35465.332
$ws.Range("I62").Value = 33263.332
$ws.Range("J62").Value = 37667.332
$ws.Range("K62").Value = 33263.332
$ws.Range("L62").Value = 37667.332
$ws.Range("M62").Value = -32639.332
$ws.Range("N62").Value = -38915.332

$ws.Range("H65").Value = 35465.332
$ws.Range("I65").Value = 33263.332
$ws.Range("J65").Value = 37667.332
$ws.Range("K65").Value = 166316.66
$ws.Range("L65").Value = 188336.66
$ws.Range("M65").Value = -163196.66
$ws.Range("N65").Value = -194576.66

$ws.Range("H126").Value = 2450.8096
$ws.Range("I126").Value = 2326.7856
$ws.Range("J126").Value = 2698.8572
$ws.Range("K126").Value = 6980.3568
$ws.Range("L126").Value = 8096.571599999999
$ws.Range("M126").Value = -4510.3568
$ws.Range("N126").Value = -13036.5716

$ws.Range("H132").Value = 3335153
$ws.Range("I132").Value = 3474089.8
$ws.Range("J132").Value = 673.5
$ws.Range("K132").Value = 10422269.4
$ws.Range("L132").Value = 2020.5
$ws.Range("M132").Value = -10419739.4
$ws.Range("N132").Value = -7080.5

$ws.Range("H136").Value = 2379706.8
$ws.Range("I136").Value = 1673818
$ws.Range("J136").Value = 3345659.8
$ws.Range("K136").Value = 5021454
$ws.Range("L136").Value = 10036979.4
$ws.Range("M136").Value = -5018904
$ws.Range("N136").Value = -10042079.4
